# Update beneficial & harmful microbiome list
# - Row 29 (last existing data row) switches from the "style 3" look to the
#   "style 4" look (same borders/font, different border-id used for the
#   previous last row of a block).
# - Ten new rows (30-39) are appended with new probiotic species, all using
#   the regular "style 3" look, except row 32 (Bifidobacterium animalis)
#   whose A/B cells get a white-fill variant of that style.
# - The active selection moves to G3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1. Row 29 restyle: copy the "style 4" look (used e.g. by row 4) onto the
#    existing row 29 cells, preserving their current values.
# ---------------------------------------------------------------------
$ws.Range("A4:D4").Copy()
$ws.Range("A29:D29").PasteSpecial($xlPasteFormats)
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2. Append new rows 30-39 using the plain "style 3" look (template: row 28)
# ---------------------------------------------------------------------
$newRows = @(
    @{ Row = 30; A = "Lactobacillus acidophilus";  B = "s__Lactobacillus_acidophilus";  D = "유익"; Height = 38 },
    @{ Row = 31; A = "Lactobacillus plantarum";    B = "s__Lactobacillus_plantarum";    D = "유익"; Height = 38 },
    @{ Row = 32; A = "Bifidobacterium animalis";   B = "s__Bifidobacterium_animalis";   D = "유익"; Height = 38 },
    @{ Row = 33; A = "Lactobacillus reuteri";      B = "s__Lactobacillus_reuteri";      D = "유익"; Height = 38 },
    @{ Row = 34; A = "Lactobacillus casei";        B = "s__Lactobacillus_casei";        D = "유익"; Height = 38 },
    @{ Row = 35; A = "Bifidobacterium breve";      B = "s__Bifidobacterium_breve";      D = "유익"; Height = 38 },
    @{ Row = 36; A = "Bifidobacterium bifidum";    B = "s__Bifidobacterium_bifidum";    D = "유익"; Height = 38 },
    @{ Row = 37; A = "Streptococcus thermophilus"; B = "s__Streptococcus_thermophilus"; D = "유익"; Height = 50.5 },
    @{ Row = 38; A = "Bifidobacterium longum";     B = "s__Bifidobacterium_longum";     D = "유익"; Height = 38 },
    @{ Row = 39; A = "Bifidobacterium longum";     B = "s__Bifidobacterium_infantis";   D = "유익"; Height = 38 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A28:D28").Copy()
    $ws.Range("A${row}:D${row}").PasteSpecial($xlPasteFormats)
    $ws.Application.CutCopyMode = $false
    $ws.Rows.Item($row).RowHeight = $r.Height

    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = $r.B
    $ws.Range("D$row").Value = $r.D
}

# Row 32 (Bifidobacterium animalis) gets a white-fill variant on A/B only
$ws.Range("A32:B32").Interior.Color = 16777215

# ---------------------------------------------------------------------
# 3. Move the active selection to G3 (matches the saved workbook state)
# ---------------------------------------------------------------------
$ws.Range("G3").Select()
